# Update PLC data 2025-10-13 13:51:09
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("LiveData")

$ws.Range("C3").Value = 163096
$ws.Range("C4").Value = 154090
$ws.Range("C7").Value = 5.52
$ws.Range("C8").Value = 64.84999999999999
